$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.406.02"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "1.805.25"
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.64"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.578"
$ws.Range("E6").Value = "  +3.98%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "35.93"
$ws.Range("E8").Value = "  +9.10%  "

$ws.Range("E9").Value = "  +2.52%  "

$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0966"
$ws.Range("E11").Value = "  +2.10%  "

$ws.Range("D12").Value = "2.065.77"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.47"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("D14").Value = "1.802.97"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("E16").Value = "  +5.48%  "

$ws.Range("D17").Value = "34.395.87"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.26"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.21"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").Value = "0.0₃0797"
$ws.Range("E20").Value = "  +0.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.51"
$ws.Range("E21").Value = "  +2.16%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("E24").Value = "  +3.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.90"
$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.91"
$ws.Range("E26").Value = "  +8.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.94"
$ws.Range("E27").Value = "  +2.63%  "

$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  +1.44%  "

$ws.Range("E31").Value = "  +1.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("E32").Value = "  +1.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("E34").Value = "  +1.09%  "

$ws.Range("D35").Value = "1.397.91"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.675"
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.48"
$ws.Range("E37").Value = "  -3.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.07"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.23"
$ws.Range("E40").Value = "  +11.56%  "

$ws.Range("E41").Value = "  +2.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "82.44"
$ws.Range("E42").Value = "  -2.34%  "

$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.51"
$ws.Range("E45").Value = "  -3.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.05"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("E47").Value = "  -5.11%  "

$ws.Range("D48").Value = "1.966.45"
$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.25"
$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("E51").Value = "  +1.42%  "
